# Update cryptocurrency price (D) and 1h volume change (E) columns
# The leading apostrophe forces Excel to store the Price value as literal
# text (matching the source data which is not a true numeric value),
# preventing values like "28.00" or "0.9990" from being normalized to 28 or 0.999.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.361.54"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "'1.935.96"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("D4").Value = "'0.9992"

$ws.Range("D5").Value = "'0.7754"
$ws.Range("E5").Value = "  +8.82%  "

$ws.Range("D6").Value = "'246.99"
$ws.Range("E6").Value = "  -1.56%  "

$ws.Range("D7").Value = "'0.9990"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "'28.00"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").Value = "'0.3215"
$ws.Range("E9").Value = "  -2.69%  "

$ws.Range("D10").Value = "'0.07097"
$ws.Range("E10").Value = "  -2.63%  "

$ws.Range("D11").Value = "'0.7839"
$ws.Range("E11").Value = "  -2.69%  "

$ws.Range("D12").Value = "'0.08022"
$ws.Range("E12").Value = "  -0.90%  "

$ws.Range("D13").Value = "'1.933.39"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "'5.380"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").Value = "'94.99"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").Value = "'14.56"
$ws.Range("E16").Value = "  -3.53%  "

$ws.Range("D17").Value = "'30.375.62"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").Value = "'255.70"
$ws.Range("E18").Value = "  +0.93%  "

$ws.Range("D19").Value = "'0.000008021"
$ws.Range("E19").Value = "  -2.04%  "

$ws.Range("D20").Value = "'5.834"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").Value = "'2.192.16"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").Value = "'0.9990"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "'0.9994"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "'6.756"
$ws.Range("E24").Value = "  -3.14%  "

$ws.Range("D25").Value = "'9.612"
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("D26").Value = "'163.91"
$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("E27").Value = "  +4.63%  "

$ws.Range("D28").Value = "'19.12"
$ws.Range("E28").Value = "  -1.01%  "

$ws.Range("D29").Value = "'2.294"
$ws.Range("E29").Value = "  -2.82%  "

$ws.Range("D30").Value = "'1.366"
$ws.Range("E30").Value = "  +1.14%  "

$ws.Range("D31").Value = "'1.523"
$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("D32").Value = "'4.431"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").Value = "'4.148"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("D34").Value = "'0.05198"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "'1.285"
$ws.Range("E35").Value = "  +1.63%  "

$ws.Range("D36").Value = "'0.7523"
$ws.Range("E36").Value = "  +0.83%  "

$ws.Range("D37").Value = "'2.772"
$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("D38").Value = "'0.01978"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").Value = "'79.11"
$ws.Range("E40").Value = "  +0.28%  "

$ws.Range("D41").Value = "'6.477"
$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").Value = "'0.4526"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").Value = "'1.985"
$ws.Range("E43").Value = "  -1.67%  "

$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "'0.8348"
$ws.Range("E45").Value = "  -1.13%  "

$ws.Range("D46").Value = "'101.39"
$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("D47").Value = "'9.850"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D48").Value = "'7.494"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("D49").Value = "'37.51"
$ws.Range("E49").Value = "  +2.08%  "

$ws.Range("D50").Value = "'982.02"
$ws.Range("E50").Value = "  +11.03%  "

$ws.Range("D51").Value = "'0.1187"
$ws.Range("E51").Value = "  +4.64%  "
